$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 287-288 (pushes existing 287-321 down to 289-323),
# matching the weekly-update pattern: new data for 2021-09-10 is prepended
# above the prior newest rows.
$ws.Rows("287:288").Insert()

# Row 287: new "Primera" quality record dated 2021-09-10 (serial 44449)
$ws.Cells.Item(287, 1).Value = 8
$ws.Cells.Item(287, 2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(287, 3).Value = 'Coquimbo'
$ws.Cells.Item(287, 4).Value = 44449
$ws.Cells.Item(287, 5).Value = 4
$ws.Cells.Item(287, 6).Value = 100112043
$ws.Cells.Item(287, 7).Value = 'Pepino ensalada'
$ws.Cells.Item(287, 8).Value = 'Sin especificar'
$ws.Cells.Item(287, 9).Value = 'Primera'
$ws.Cells.Item(287, 10).Value = 800
$ws.Cells.Item(287, 11).Value = 15000
$ws.Cells.Item(287, 12).Value = 16000
$ws.Cells.Item(287, 13).Value = 15500
$ws.Cells.Item(287, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(287, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(287, 16).Value = 258
$ws.Cells.Item(287, 17).Value = 60
$ws.Cells.Item(287, 18).Value = 'Hortaliza'

# Row 288: new "Segunda" quality record dated 2021-09-10 (serial 44449)
$ws.Cells.Item(288, 1).Value = 8
$ws.Cells.Item(288, 2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(288, 3).Value = 'Coquimbo'
$ws.Cells.Item(288, 4).Value = 44449
$ws.Cells.Item(288, 5).Value = 4
$ws.Cells.Item(288, 6).Value = 100112043
$ws.Cells.Item(288, 7).Value = 'Pepino ensalada'
$ws.Cells.Item(288, 8).Value = 'Sin especificar'
$ws.Cells.Item(288, 9).Value = 'Segunda'
$ws.Cells.Item(288, 10).Value = 560
$ws.Cells.Item(288, 11).Value = 10000
$ws.Cells.Item(288, 12).Value = 11000
$ws.Cells.Item(288, 13).Value = 10500
$ws.Cells.Item(288, 14).Value = '$/caja 100 unidades'
$ws.Cells.Item(288, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(288, 16).Value = 105
$ws.Cells.Item(288, 17).Value = 100
$ws.Cells.Item(288, 18).Value = 'Hortaliza'
